$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph entirely. It currently reads:
#      "Meta description" (bold) + ": Read our review of Dia De Los
#      Muertos, a colorful and unique slot game celebrating the Mexican
#      tradition of Dia de Los Muertos. Play for free now!"
# ---------------------------------------------------------------------------
$metaFind = $d.Content.Duplicate
$metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFind.Find.Found) {
    # Resolve which paragraph the match falls in, then delete that whole
    # paragraph (Range.Paragraphs is unreliable in this host, so compute the
    # index from how many paragraph marks precede the match instead).
    $precedingRange = $d.Range(0, $metaFind.Start)
    $metaParaIndex = $precedingRange.Paragraphs.Count + 1
    $metaPara = $d.Paragraphs.Item($metaParaIndex)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. At the end of the document, insert a new bold paragraph with the page
#    title right before the final paragraph, and change the final
#    paragraph's (italic) text from the old image-generation prompt to the
#    review blurb that used to live in the meta-description paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End

# Remove the existing (italic) text of the last paragraph, keeping the
# paragraph mark itself so the section properties stay attached correctly.
$lastContent = $d.Range($lastStart, $lastEnd)
$lastContent.Delete()

# Rebuild both paragraphs from scratch via raw OOXML so the run/formatting
# layout matches exactly: a leading empty run followed by a single styled
# text run, with no stray paragraph-mark formatting left behind.
$insRange = $d.Range($lastStart, $lastStart)
$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dia De Los Muertos Free: Colorful and Unique Slot Game</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Dia De Los Muertos, a colorful and unique slot game celebrating the Mexican tradition of Dia de Los Muertos. Play for free now!</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insRange.InsertXML($xmlSnippet)
